# Insert a new record row at row 97 (pushes the existing rows 97-142 down
# to 98-143, matching the dimension change A1:R142 -> A1:R143), then
# populate the new row with the new "Pepino ensalada" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 97..142 down to 98..143, inserting a blank row at 97.
$ws.Rows(97).Insert()

# Fill in the new row 97 with the new record.
$ws.Cells.Item(97, 1).Value = 11
$ws.Cells.Item(97, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(97, 3).Value = "Bíobío"
$ws.Cells.Item(97, 4).Value = 44755
$ws.Cells.Item(97, 5).Value = 8
$ws.Cells.Item(97, 6).Value = 100112043
$ws.Cells.Item(97, 7).Value = "Pepino ensalada"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 100
$ws.Cells.Item(97, 11).Value = 19000
$ws.Cells.Item(97, 12).Value = 20000
$ws.Cells.Item(97, 13).Value = 19500
$ws.Cells.Item(97, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(97, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(97, 16).Value = 325
$ws.Cells.Item(97, 17).Value = 60
$ws.Cells.Item(97, 18).Value = "Hortaliza"
